$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 59, shifting existing rows 59:140 down to 60:140
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new weekly price record
$ws.Cells.Item(59, 1).Value = 4
$ws.Cells.Item(59, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(59, 3).Value = "Los Lagos"
$ws.Cells.Item(59, 4).Value = 44467
$ws.Cells.Item(59, 5).Value = 10
$ws.Cells.Item(59, 6).Value = "Fruta"
$ws.Cells.Item(59, 7).Value = 100108
$ws.Cells.Item(59, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(59, 9).Value = 100108005
$ws.Cells.Item(59, 10).Value = "Piña"
$ws.Cells.Item(59, 11).Value = "Caramelo"
$ws.Cells.Item(59, 12).Value = "Primera"
$ws.Cells.Item(59, 13).Value = 160
$ws.Cells.Item(59, 14).Value = 22000
$ws.Cells.Item(59, 15).Value = 22000
$ws.Cells.Item(59, 16).Value = 22000
$ws.Cells.Item(59, 17).Value = "`$/caja 12 unidades"
$ws.Cells.Item(59, 18).Value = "Ecuador"
$ws.Cells.Item(59, 19).Value = 1833
$ws.Cells.Item(59, 20).Value = 12
